$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    '126 Racecourse Road Public Housing Tower Flemington',
    '139 Highett St Apartment Complex Richmond',
    '3175 The Bays Aged Care Facility Hastings',
    '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
    'Al-Taqwa College Truganina',
    'Allbright Manor Aged Care Croydon North Tier 1A',
    'Australia Post Distribution Centre Sunshine West',
    'Australian Lamb Colac East',
    'Baker Bleu Caulfield North',
    'Baxter Foods Australia Campbellfield',
    'CS Square Caroline Springs',
    'Cafe Roco Dandenong',
    'Campbellfield Ford Complex Vaccination Clinic Campbellfield',
    'Cardinia Lakes Early Learning Centre Pakenham',
    'Carton Finishing Pty. Ltd. Campbellfield',
    'Chemist Warehouse Fillo Drive Somerton',
    'Coles Coburg North Village',
    'Coles Pakenham Place Shopping Centre',
    'Community Kids Bayswater Early Education Centre Bayswater North',
    'Construction Site Olea Apartment Caulfield North',
    'Costco Wholesale Epping',
    'Crusader Caravans Epping',
    'Dandenong Police Station Dandenong',
    'DayHab Rehabilitation Treatment Centre Ringwood East',
    'Embracia Aged Care Reservoir Outbreak',
    'Ermha365 Residential Disability Care Services Doveton',
    'FedEx Station Melbourne Airport',
    'Fine Food Holdings Pty Ltd Dandenong South',
    'Fonterra Manufacturing Workplace Campbellfield',
    'General Foods Campbellfield',
    'Gladstone Parade Early Learning & Kinder Glenroy',
    'Goodstart Early Learning Altona',
    'Green Leaves Early Learning Centre Highlands Craigieburn',
    'Greenvale Primary School',
    'HEI Schools Emerald Early Learning Centre Emerald',
    'Hello Fresh Warehouse Ravenhall',
    'Ibis Kingsgate Hotel Melbourne',
    'Industrial Galvanizers Valmont Coatings Campbellfield',
    'Inghams Enterprises Thomastown',
    'Kool Kidz Childcare Narre Warren',
    'Lantmannen Unibake Australia Mordialloc',
    'Linfox Somerton National Distribution Centre Somerton',
    'Mecca Distribution Centre Warehouse Melbourne Airport',
    'Melbourne Assessment Prison West Melbourne',
    'Melbourne Metropolitan Remand Centre Ravenhall',
    'Melbourne West Police Station Docklands',
    'Mill Park Police Station Mill Park',
    'MyCentre Childcare Broadmeadows',
    'Nido Early School Ascot Vale',
    'Nido Early School Glenroy',
    'Northern Health Northern Hospital Epping Emergency Department Tier 1B',
    'Northern Health The Northern Hospital Epping',
    'OnQ Plumbing and Excavations Craigieburn',
    'Pacific Meat Thomastown',
    'Private Residence Daycare Allumba Way Wollert',
    'Ravenhall Correctional Centre Ravenhall',
    'Richmond Quarter 261-271 Bridge Road Construction Site Richmond',
    'St Margaret''s Primary School OSHC Maribyrnong',
    'St Vincents Hospital Emergency Department Melbourne',
    'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
    'The Royal Melbourne Hospital Parkville',
    'The Royal Melbourne Hospital Parkville Emergency Department',
    'The Royal Melbourne Hospital Ward 6SE Parkville',
    'The Royal Talbot Rehabilitation Centre Kew',
    'ThorwestenCabinets Pakenham',
    'Truganina Early Learning Centre Truganina',
    'Visy Recycling Springvale',
    'Wallaby Childcare Wollert',
    'Werribee Mercy Hospital Emergency Department',
    'Western Health Footscray Hospital Emergency Department',
    'Western Health Sunshine Hospital Emergency Department',
    'Western Health Sunshine Hospital GEM Ward St Albans'
)

$values = @(6, 10, 9, 12, 6, 8, 7, 10, 7, 5, 11, 6, 9, 6, 12, 7, 6, 5, 17, 6, 16, 12, 9, 7, 21, 10, 15, 10, 8, 12, 7, 6, 16, 5, 5, 7, 6, 5, 5, 15, 20, 6, 7, 9, 9, 8, 9, 7, 28, 23, 53, 23, 7, 5, 8, 10, 5, 12, 20, 17, 6, 6, 9, 11, 8, 8, 29, 17, 15, 8, 15, 6)

$rowCount = $names.Length

$data = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $names[$i]
    $data[$i,1] = $values[$i]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$writeAddr = "A" + $startRow + ":B" + $endRow
$ws.Range($writeAddr).Value = $data

# The source sheet previously had data through row 81; the refreshed
# data only spans through row $endRow, so clear out the old trailing rows.
$lastRow = 81
if ($lastRow -gt $endRow) {
    $clearAddr = "A" + ($endRow + 1) + ":B" + $lastRow
    $ws.Range($clearAddr).ClearContents()
}
